$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (activity counts) and column C (hour) values for rows 2-11
$ws.Range("B2").Value = 2432
$ws.Range("C2").Value = 17

$ws.Range("B3").Value = 216
$ws.Range("C3").Value = 17

$ws.Range("B4").Value = 510
$ws.Range("C4").Value = 17

$ws.Range("B5").Value = 1434
$ws.Range("C5").Value = 17

$ws.Range("B6").Value = 1232
$ws.Range("C6").Value = 17

$ws.Range("B7").Value = 324
$ws.Range("C7").Value = 17

$ws.Range("B8").Value = 1013
$ws.Range("C8").Value = 17

$ws.Range("B9").Value = 164
$ws.Range("C9").Value = 17

$ws.Range("C10").Value = 17

$ws.Range("B11").Value = 111
$ws.Range("C11").Value = 17
